# Modified as per new Assessment Creation flow
# Updates the CreateAssessment sheet: fixes a StartTime value, shrinks two
# instruction-row heights, and appends a new "CreateFree" assessment row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CreateAssessment")

# --- StartTime fix on the Group assessment row (row 3) ---
$ws.Range("I3").Value = "'30/03/20 17:20"

# --- Row height adjustments for the Mock / Group instruction rows ---
$ws.Rows.Item(2).RowHeight = 87
$ws.Rows.Item(3).RowHeight = 125

# --- New row 6: CreateFree / Free_Assessment test case ---
$ws.Range("A6").Value = 5
$ws.Range("B6").Value = "CreateFree"
$ws.Range("C6").Value = "Free_Assessment"
$ws.Range("D6").Value = "Free_Assessment_Description"

$freeInstructions = "1.This online test requires a stable internet connection.`n2.We recommend using Chrome or Firefox browsers only.`n3.This test must be taken on the desktop only.`n4.The total time to complete the test is 60 mins(1 hr).`n5.This is a free assessment.`n7.If you run out of time, all attempted questions will be auto-submitted.`n8.For any technical queries, please email us at  support@knowledgehut.com."
$ws.Range("E6").Value = $freeInstructions
$ws.Range("E6").WrapText = $true

$ws.Range("F6").Value = "'100"
$ws.Range("G6").Value = "'30"
$ws.Range("H6").Value = "'"
$ws.Range("I6").Value = "'"
$ws.Range("J6").Value = "'"
$ws.Range("K6").Value = "'20"

$ws.Rows.Item(6).RowHeight = 135

# --- Selection / view state ---
$ws.Activate()
$ws.Range("E6").Select()
